# Fruta / hortaliza, semanal
# Insert two new weekly records at the top of the data block (rows 1078-1174),
# pushing the existing data down by two rows (new rows 1080-1176).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 1078 (existing row 1078 and everything below it
# shifts down by one each time this is called).
$ws.Rows.Item(1078).Insert()
$ws.Rows.Item(1078).Insert()

# --- New row 1078: Packham's Triumph, Primera ---
$ws.Cells.Item(1078, 1).Value = 8
$ws.Cells.Item(1078, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1078, 3).Value = "Coquimbo"
$ws.Cells.Item(1078, 4).Value = 45212
$ws.Cells.Item(1078, 5).Value = 4
$ws.Cells.Item(1078, 6).Value = "Fruta"
$ws.Cells.Item(1078, 7).Value = 100104
$ws.Cells.Item(1078, 8).Value = "Frutos de pepita"
$ws.Cells.Item(1078, 9).Value = 100104005
$ws.Cells.Item(1078, 10).Value = "Pera"
$ws.Cells.Item(1078, 11).Value = "Packham's Triumph"
$ws.Cells.Item(1078, 12).Value = "Primera"
$ws.Cells.Item(1078, 13).Value = 10
$ws.Cells.Item(1078, 14).Value = 370000
$ws.Cells.Item(1078, 15).Value = 380000
$ws.Cells.Item(1078, 16).Value = 375000
$ws.Cells.Item(1078, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(1078, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(1078, 19).Value = 833
$ws.Cells.Item(1078, 20).Value = 450
$ws.Range("D1078").NumberFormat = $ws.Range("D1080").NumberFormat

# --- New row 1079: Packham's Triumph, Segunda ---
$ws.Cells.Item(1079, 1).Value = 8
$ws.Cells.Item(1079, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1079, 3).Value = "Coquimbo"
$ws.Cells.Item(1079, 4).Value = 45212
$ws.Cells.Item(1079, 5).Value = 4
$ws.Cells.Item(1079, 6).Value = "Fruta"
$ws.Cells.Item(1079, 7).Value = 100104
$ws.Cells.Item(1079, 8).Value = "Frutos de pepita"
$ws.Cells.Item(1079, 9).Value = 100104005
$ws.Cells.Item(1079, 10).Value = "Pera"
$ws.Cells.Item(1079, 11).Value = "Packham's Triumph"
$ws.Cells.Item(1079, 12).Value = "Segunda"
$ws.Cells.Item(1079, 13).Value = 10
$ws.Cells.Item(1079, 14).Value = 340000
$ws.Cells.Item(1079, 15).Value = 350000
$ws.Cells.Item(1079, 16).Value = 345000
$ws.Cells.Item(1079, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(1079, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(1079, 19).Value = 767
$ws.Cells.Item(1079, 20).Value = 450
$ws.Range("D1079").NumberFormat = $ws.Range("D1080").NumberFormat
